# Atualiza o relatorio PROMAR com os novos dados (filtro por idade/faixa etaria)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insere as novas linhas necessarias para acomodar os dados adicionais ---
# Inserir de baixo para cima para que os numeros de linha originais continuem validos.
$ws.Rows(30).Insert()
$ws.Rows(25).Insert()
$ws.Rows(14).Insert()

# --- Secao "Categoria" / "Valor" (linhas 1-4) ---
$ws.Range("B2").Value = 8
$ws.Range("B4").Value = 5

# --- Secao "Respostas Mais Acertadas" (linhas 6-13) ---
$ws.Range("B8").Value = 4
$ws.Range("A9").Value = "b-Um ser humano"
$ws.Range("B9").Value = 2
$ws.Range("A10").Value = "b-Brincando na praia"
$ws.Range("B10").Value = 2
$ws.Range("A11").Value = "a-Um animal"
$ws.Range("B11").Value = 2
$ws.Range("A12").Value = "c-Olhando os peixes"
$ws.Range("B12").Value = 1
$ws.Range("A13").Value = "a-Não jogando lixo nele"
$ws.Range("B13").Value = 1

# --- Secao "Cidades com Melhor Desempenho" (linhas 15-18) ---
$ws.Range("B17").Value = 6
$ws.Range("B18").Value = 1

# --- Secao "Notas Mais Dadas" (linhas 20-25) ---
$ws.Range("B22").Value = 5
$ws.Range("A24").Value = 5
$ws.Range("B24").Value = 1
$ws.Range("A25").Value = 2
$ws.Range("B25").Value = 1

# --- Secao "Idades Mais Visitadas" (linhas 27-32) ---
$ws.Range("A30").Value = "5-7"
$ws.Range("B30").Value = 2
$ws.Range("A31").Value = "mais_de_12"
$ws.Range("B31").Value = 1
$ws.Range("A32").Value = "7-9"
$ws.Range("B32").Value = 1
